$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3552769422531128
$ws.Range("B1").Value = 0.9985016584396362
$ws.Range("C1").Value = 4.787903785705566
$ws.Range("D1").Value = 1.998209118843079
$ws.Range("E1").Value = 0.8601756691932678
